# Update factsheets with text edits from COMM
# Converts numeric "count" cells to text cells (matching formatted numbers
# with thousands separators where applicable), and adds a "Total" summary
# row to the County sheet.

$wb = $excel.ActiveWorkbook

function Set-TextValue($Range, $Text) {
    $Range.NumberFormat = "@"
    $Range.Value = $Text
}

# ----- Sheet "Overall" -----
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Range("A2") "1,592"

# ----- Sheet "County" -----
$wsCounty = $wb.Worksheets.Item("County")
Set-TextValue $wsCounty.Range("B2") "1,592"

# Add new "Total" row (row 3) mirroring the existing data row.
Set-TextValue $wsCounty.Range("A3") "Total"
Set-TextValue $wsCounty.Range("B3") "1,592"
Set-TextValue $wsCounty.Range("C3") "$18,688,408,098"
Set-TextValue $wsCounty.Range("D3") "9.72%"
Set-TextValue $wsCounty.Range("E3") "-9.87%"
Set-TextValue $wsCounty.Range("F3") "62.94%"

# ----- Sheet "Congressional District" -----
$wsCong = $wb.Worksheets.Item("Congressional District")
Set-TextValue $wsCong.Range("B2") "1,592"
Set-TextValue $wsCong.Range("B3") "1,592"

# ----- Sheet "Size" -----
$wsSize = $wb.Worksheets.Item("Size")
Set-TextValue $wsSize.Range("B2") "226"
Set-TextValue $wsSize.Range("B3") "594"
Set-TextValue $wsSize.Range("B4") "235"
Set-TextValue $wsSize.Range("B5") "191"
Set-TextValue $wsSize.Range("B6") "329"
Set-TextValue $wsSize.Range("B7") "17"
Set-TextValue $wsSize.Range("B8") "1,592"

# ----- Sheet "Subsector" -----
$wsSub = $wb.Worksheets.Item("Subsector")
Set-TextValue $wsSub.Range("B2") "155"
Set-TextValue $wsSub.Range("B3") "212"
Set-TextValue $wsSub.Range("B4") "77"
Set-TextValue $wsSub.Range("B5") "111"
Set-TextValue $wsSub.Range("B6") "3"
Set-TextValue $wsSub.Range("B7") "315"
Set-TextValue $wsSub.Range("B8") "142"
Set-TextValue $wsSub.Range("B9") "2"
Set-TextValue $wsSub.Range("B10") "241"
Set-TextValue $wsSub.Range("B11") "22"
Set-TextValue $wsSub.Range("B12") "304"
Set-TextValue $wsSub.Range("B13") "8"
Set-TextValue $wsSub.Range("B14") "1,592"
